$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'62.403.11"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -1.99%  "

$cell = $ws.Range("D3")
$cell.Value = "'3.162.34"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -3.64%  "

$cell = $ws.Range("D5")
$cell.Value = "'585.80"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -3.04%  "

$cell = $ws.Range("D6")
$cell.Value = "'134.80"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -4.71%  "

$ws.Range("E7").Value = "  -0.14%  "

$cell = $ws.Range("D8")
$cell.Value = "'3.160.71"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -3.60%  "

$cell = $ws.Range("D9")
$cell.Value = "'0.506"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -2.23%  "

$ws.Range("E10").Value = "  -5.37%  "

$ws.Range("E11").Value = "  -2.93%  "

$cell = $ws.Range("D12")
$cell.Value = "'0.453"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -3.19%  "

$ws.Range("E13").Value = "  -4.86%  "

$cell = $ws.Range("D14")
$cell.Value = "'33.26"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.44%  "

$cell = $ws.Range("D15")
$cell.Value = "'3.685.82"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -3.63%  "

$ws.Range("E16").Value = "  -2.04%  "

$cell = $ws.Range("D17")
$cell.Value = "'3.169.02"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -3.39%  "

$cell = $ws.Range("D18")
$cell.Value = "'62.391.51"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -2.09%  "

$cell = $ws.Range("D19")
$cell.Value = "'6.51"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -4.53%  "

$cell = $ws.Range("D20")
$cell.Value = "'454.35"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -5.21%  "

$cell = $ws.Range("D21")
$cell.Value = "'13.91"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.86%  "

$cell = $ws.Range("D22")
$cell.Value = "'0.699"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -3.79%  "

$cell = $ws.Range("D23")
$cell.Value = "'7.59"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -4.86%  "

$cell = $ws.Range("D24")
$cell.Value = "'83.49"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.81%  "

$cell = $ws.Range("D25")
$cell.Value = "'13.23"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.02%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("B27").Value = "FirstDigitalUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Range("D27")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell = $ws.Range("D28")
$cell.Value = "'2.68"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -3.04%  "

$cell = $ws.Range("D29")
$cell.Value = "'6.82"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -6.44%  "

$cell = $ws.Range("D30")
$cell.Value = "'7.71"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -4.73%  "

$ws.Range("E31").Value = "  -7.01%  "

$cell = $ws.Range("D32")
$cell.Value = "'27.12"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -5.41%  "

$ws.Range("E33").Value = "  -1.48%  "

$cell = $ws.Range("D34")
$cell.Value = "'2.37"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -5.98%  "

$ws.Range("E35").Value = "  -5.98%  "

$cell = $ws.Range("D36")
$cell.Value = "'5.90"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.65%  "

$cell = $ws.Range("D37")
$cell.Value = "'51.14"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.45%  "

$cell = $ws.Range("D38")
$cell.Value = "'0.0₃0693"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -6.09%  "

$cell = $ws.Range("D39")
$cell.Value = "'0.0382"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -3.72%  "

$cell = $ws.Range("D40")
$cell.Value = "'2.72"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "

$cell = $ws.Range("D41")
$cell.Value = "'395.28"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -7.18%  "

$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$cell = $ws.Range("D42")
$cell.Value = "'7.98"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -4.08%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D43")
$cell.Value = "'0.112"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.45%  "

$cell = $ws.Range("D44")
$cell.Value = "'2.794.50"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -8.44%  "

$cell = $ws.Range("D45")
$cell.Value = "'0.249"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -5.44%  "

$cell = $ws.Range("D47")
$cell.Value = "'2.12"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.40%  "

$cell = $ws.Range("D48")
$cell.Value = "'35.47"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +2.07%  "

$cell = $ws.Range("D49")
$cell.Value = "'125.09"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.31%  "

$cell = $ws.Range("D50")
$cell.Value = "'25.23"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.42%  "

$ws.Range("E51").Value = "  -3.84%  "
